$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the current row 117,
# pushing the existing row 117 (and everything after it) down by one row.
$ws.Rows(117).Insert()

$ws.Range("A117").Value = 10
$ws.Range("B117").Value = "Vega Modelo de Temuco"
$ws.Range("C117").Value = "La Araucanía"
$ws.Range("D117").Value = 44447
$ws.Range("E117").Value = 9
$ws.Range("F117").Value = 100112040
$ws.Range("G117").Value = "Cilantro"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 65
$ws.Range("K117").Value = 4000
$ws.Range("L117").Value = 4500
$ws.Range("M117").Value = 4269
$ws.Range("N117").Value = "$/docena de atados (2 kilos)"
$ws.Range("O117").Value = "Región Metropolitana"
$ws.Range("P117").Value = 2134
$ws.Range("Q117").Value = 2
$ws.Range("R117").Value = "Hortaliza"
